$d = $word.ActiveDocument

$pairs = @(
    @("393÷9=", "862÷4="),
    @("810÷2=", "510÷6="),
    @("704÷2=", "868÷3="),
    @("120÷6=", "118÷8="),
    @("771÷2=", "318÷6="),
    @("545÷6=", "700÷3="),
    @("317÷2=", "789÷7="),
    @("886÷7=", "944÷7="),
    @("220÷9=", "889÷7="),
    @("858÷2=", "957÷8="),
    @("217÷6=", "320÷8="),
    @("216÷9=", "609÷7="),
    @("648÷7=", "258÷4="),
    @("467÷2=", "131÷8="),
    @("761÷4=", "878÷4="),
    @("263÷7=", "313÷5="),
    @("743÷9=", "701÷6="),
    @("319÷6=", "298÷2="),
    @("179÷4=", "106÷2="),
    @("286÷8=", "651÷7="),
    @("139÷2=", "999÷2="),
    @("537÷3=", "354÷3="),
    @("447÷3=", "271÷5="),
    @("164÷6=", "290÷7="),
    @("755÷5=", "511÷9=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
